$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

$ws.Range("B2").Value = "yes"
$ws.Range("B3").Value = "yes"
$ws.Range("B5").Value = "yes"
$ws.Range("B6").Value = "yes"
$ws.Range("B7").Value = "yes"

$ws.Range("B2:B7").Select()
